# Auto-generated script applying scheduled market-price refresh to Sheets/Masamune_Profits.xlsx
# Updates computed price/profit columns (H-N) on the ALC/ARM/BSM/CRP/CUL/GSM/LTW worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1039.875
$ws.Range("I18").Value = 869.8333
$ws.Range("J18").Value = 1550
$ws.Range("K18").Value = 869.8333
$ws.Range("L18").Value = 1550
$ws.Range("M18").Value = -585.8333
$ws.Range("N18").Value = -2118
$ws.Range("H31").Value = 1499.5
$ws.Range("I31").Value = 1499.5
$ws.Range("K31").Value = 4498.5
$ws.Range("M31").Value = -4268.5
$ws.Range("H40").Value = 11940.1
$ws.Range("I40").Value = 16080.143
$ws.Range("J40").Value = 2280
$ws.Range("K40").Value = 16080.143
$ws.Range("L40").Value = 2280
$ws.Range("M40").Value = -15905.143
$ws.Range("N40").Value = -2630
$ws.Range("H98").Value = 54155.8
$ws.Range("J98").Value = 88680.78
$ws.Range("L98").Value = 88680.78
$ws.Range("N98").Value = -91676.78
$ws.Range("H106").Value = 196665.5
$ws.Range("I106").Value = 2905
$ws.Range("J106").Value = 261252.33
$ws.Range("K106").Value = 2905
$ws.Range("L106").Value = 261252.33
$ws.Range("M106").Value = -2274
$ws.Range("N106").Value = -262514.33
$ws.Range("H122").Value = 54155.8
$ws.Range("J122").Value = 88680.78
$ws.Range("L122").Value = 266042.34
$ws.Range("N122").Value = -270942.34
$ws.Range("H137").Value = 1430214
$ws.Range("I137").Value = 2083674.8
$ws.Range("J137").Value = 7975.9414
$ws.Range("K137").Value = 6251024.4
$ws.Range("L137").Value = 23927.8242
$ws.Range("M137").Value = -6248474.4
$ws.Range("N137").Value = -29027.8242
$ws.Range("H138").Value = 2043.527
$ws.Range("I138").Value = 1655.9584
$ws.Range("J138").Value = 2229.56
$ws.Range("K138").Value = 4967.8752
$ws.Range("L138").Value = 6688.68
$ws.Range("M138").Value = 172.1247999999996
$ws.Range("N138").Value = -16968.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1844.7142
$ws.Range("I74").Value = 1480.5641
$ws.Range("J74").Value = 3264.9
$ws.Range("K74").Value = 1480.5641
$ws.Range("L74").Value = 3264.9
$ws.Range("M74").Value = -606.5641000000001
$ws.Range("N74").Value = -5012.9
$ws.Range("H77").Value = 1844.7142
$ws.Range("I77").Value = 1480.5641
$ws.Range("J77").Value = 3264.9
$ws.Range("K77").Value = 7402.8205
$ws.Range("L77").Value = 16324.5
$ws.Range("M77").Value = -3034.8205
$ws.Range("N77").Value = -25060.5
$ws.Range("H110").Value = 1461.7778
$ws.Range("I110").Value = 1458.9131
$ws.Range("J110").Value = 1478.25
$ws.Range("K110").Value = 1458.9131
$ws.Range("L110").Value = 1478.25
$ws.Range("M110").Value = 586.0869
$ws.Range("N110").Value = -5568.25
$ws.Range("H122").Value = 1554.4348
$ws.Range("I122").Value = 1316
$ws.Range("J122").Value = 1925.3334
$ws.Range("K122").Value = 3948
$ws.Range("L122").Value = 5776.0002
$ws.Range("M122").Value = -1498
$ws.Range("N122").Value = -10676.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2386.25
$ws.Range("I99").Value = 2279
$ws.Range("J99").Value = 2565
$ws.Range("K99").Value = 2279
$ws.Range("L99").Value = 2565
$ws.Range("M99").Value = -781
$ws.Range("N99").Value = -5561
$ws.Range("H105").Value = 2153.6428
$ws.Range("I105").Value = 2047.8148
$ws.Range("K105").Value = 2047.8148
$ws.Range("M105").Value = -300.8148000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 70003
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("H6").Value = 3983.3333
$ws.Range("I6").Value = 3666.6667
$ws.Range("J6").Value = 4300
$ws.Range("K6").Value = 3666.6667
$ws.Range("L6").Value = 4300
$ws.Range("M6").Value = -3553.6667
$ws.Range("N6").Value = -4526
$ws.Range("H7").Value = 1353.5
$ws.Range("I7").Value = 1154.5625
$ws.Range("J7").Value = 2945
$ws.Range("K7").Value = 1154.5625
$ws.Range("L7").Value = 2945
$ws.Range("M7").Value = -1041.5625
$ws.Range("N7").Value = -3171
$ws.Range("H17").Value = 39990
$ws.Range("J17").Value = 39990
$ws.Range("L17").Value = 39990
$ws.Range("N17").Value = -40338
$ws.Range("H25").Value = 33555.6
$ws.Range("I25").Value = 9800
$ws.Range("J25").Value = 39494.5
$ws.Range("K25").Value = 9800
$ws.Range("L25").Value = 39494.5
$ws.Range("M25").Value = -9626
$ws.Range("N25").Value = -39842.5
$ws.Range("H50").Value = 32149.8
$ws.Range("J50").Value = 38937.25
$ws.Range("L50").Value = 38937.25
$ws.Range("N50").Value = -40187.25
$ws.Range("H51").Value = 100024450
$ws.Range("I51").Value = 250002940
$ws.Range("J51").Value = 38773.332
$ws.Range("K51").Value = 250002940
$ws.Range("L51").Value = 38773.332
$ws.Range("M51").Value = -250002204
$ws.Range("N51").Value = -40245.332
$ws.Range("H59").Value = 30720
$ws.Range("J59").Value = 30720
$ws.Range("L59").Value = 30720
$ws.Range("N59").Value = -33010
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = ""
$ws.Range("H61").Value = 100024450
$ws.Range("I61").Value = 250002940
$ws.Range("J61").Value = 38773.332
$ws.Range("K61").Value = 250002940
$ws.Range("L61").Value = 38773.332
$ws.Range("M61").Value = -250002592
$ws.Range("N61").Value = -39469.332
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H96").Value = 47909.75
$ws.Range("J96").Value = 47909.75
$ws.Range("L96").Value = 47909.75
$ws.Range("N96").Value = -53401.75
$ws.Range("H125").Value = 19666.666
$ws.Range("J125").Value = 19666.666
$ws.Range("L125").Value = 19666.666
$ws.Range("N125").Value = -24586.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1236.0244
$ws.Range("J68").Value = 1295.3016
$ws.Range("L68").Value = 3885.9048
$ws.Range("N68").Value = -5507.9048
$ws.Range("H71").Value = 1236.0244
$ws.Range("J71").Value = 1295.3016
$ws.Range("L71").Value = 11657.7144
$ws.Range("N71").Value = -19769.7144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 251.5
$ws.Range("I4").Value = 251.5
$ws.Range("K4").Value = 251.5
$ws.Range("M4").Value = -139.5
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 23575.455
$ws.Range("I94").Value = 10000
$ws.Range("J94").Value = 26592.223
$ws.Range("K94").Value = 10000
$ws.Range("L94").Value = 26592.223
$ws.Range("N94").Value = -27944.223
